# The league's results were re-matched to the correct fixtures: the data
# (bookmaker id, teams, score, odds, ...) that had been attached to one row
# actually belongs to the other row of each pair, and vice versa.
# Column A (sequential match number) together with columns C/D/E (Div, Div
# Original Name, Date) stay put; every other column (B, F..AC) is swapped
# between the two rows of each pair.

function Swap-Cell([string]$col, [int]$row1, [int]$row2) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value()
    $v2 = $ws.Range($addr2).Value()
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}

function Swap-Rows([int]$row1, [int]$row2) {
    Swap-Cell "B" $row1 $row2
    Swap-Cell "F" $row1 $row2
    Swap-Cell "G" $row1 $row2
    Swap-Cell "H" $row1 $row2
    Swap-Cell "I" $row1 $row2
    Swap-Cell "J" $row1 $row2
    Swap-Cell "K" $row1 $row2
    Swap-Cell "L" $row1 $row2
    Swap-Cell "M" $row1 $row2
    Swap-Cell "N" $row1 $row2
    Swap-Cell "O" $row1 $row2
    Swap-Cell "P" $row1 $row2
    Swap-Cell "Q" $row1 $row2
    Swap-Cell "R" $row1 $row2
    Swap-Cell "S" $row1 $row2
    Swap-Cell "T" $row1 $row2
    Swap-Cell "U" $row1 $row2
    Swap-Cell "V" $row1 $row2
    Swap-Cell "W" $row1 $row2
    Swap-Cell "X" $row1 $row2
    Swap-Cell "Y" $row1 $row2
    Swap-Cell "Z" $row1 $row2
    Swap-Cell "AA" $row1 $row2
    Swap-Cell "AB" $row1 $row2
    Swap-Cell "AC" $row1 $row2
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Swap-Rows 39 40
Swap-Rows 111 112
